$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing header labels ("DigitalInstantition" -> "DigitalInstantiation")
$ws.Range("F1").Value = "DigitalInstantiation.aapb_preservation_lto"
$ws.Range("G1").Value = "DigitalInstantiation.aapb_preservation_disk"

# Add new column H: md5 header + test value
$ws.Range("H1").Value = "DigitalInstantiation.md5"
$ws.Range("H2").Value = "qwertyqwerty"

# Match the selection left behind by the author's edit
$ws.Range("G1").Select()
